# Generate Report for Handback
# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps for the 22517586-... rows on the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "2016-03-21 06:37:24"
$wsZhCn.Range("H4").Value = "2016-03-21 06:37:52"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "2016-03-21 06:37:27"
$wsDeDe.Range("H4").Value = "2016-03-21 06:37:58"
